$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C1").Value = "Berg"
$ws.Range("D1").Value = "Nordosten"
$ws.Range("E1").Value = "Süd"
